# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (fund-holding detail) right after "2021-Q4"
#   and before "总计", populated like the other quarterly detail sheets.
# - Insert a new top data row in "总计" for 2022-Q1, pushing the existing
#   quarters down by one row and renumbering the index column.

$wb = $excel.ActiveWorkbook

function Set-TextValue($sheet, $row, $col, $text) {
    # Force the cell to be stored as text (not auto-coerced to a number),
    # even when the value looks numeric (e.g. "159869" or "6.20").
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" detail sheet by duplicating "2021-Q4" (same
#    columns/header/style) and overwriting its data with the new figures.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$fundRows = @(
    @{ Code = "159869"; Name = "华夏中证动漫游戏ETF";   Size = "6.20"; Pos = "98.75"; Ratio = "3.99"; Value = "0.2474"; Rank = 8 },
    @{ Code = "516010"; Name = "国泰中证动漫游戏ETF";   Size = "4.95"; Pos = "98.91"; Ratio = "3.90"; Value = "0.1930"; Rank = 8 },
    @{ Code = "516770"; Name = "华泰柏瑞中证动漫游戏ETF"; Size = "1.11"; Pos = "96.56"; Ratio = "3.94"; Value = "0.0437"; Rank = 8 }
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    Set-TextValue $newSheet $r 2 $row.Code
    Set-TextValue $newSheet $r 3 $row.Name
    Set-TextValue $newSheet $r 4 $row.Size
    Set-TextValue $newSheet $r 5 $row.Pos
    Set-TextValue $newSheet $r 6 $row.Ratio
    Set-TextValue $newSheet $r 7 $row.Value
    $newSheet.Cells.Item($r, 8).Value = $row.Rank
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row right under the
#    header for 2022-Q1 and renumber the existing index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.48

for ($r = 3; $r -le 6; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
